# Apply updated K column (strikeouts) values: G2:G67
# These were regenerated upstream (commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals") and are written here as the
# new literal s_vals for column G ("K").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,0,0,0,1,2,1,1,3,0,3,0,1,2,1,1,1,1,1,0,0,1,1,1,0,1,1,2,0,1,2,3,0,1,0,1,3,1,2,1,0,1,1,2,0,3,1,2,0,0,0,1,1,1,0,0,2,0,1,1,1,0,0,1,1,2)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
